$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1: "time_taken", formatted like the other header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# time_taken values for each data row (plain text, matching microsecond-precision source)
$ws.Range("F2").Value = "2021-10-05 13:38:29.769959"
$ws.Range("F3").Value = "2021-10-05 13:38:29.769970"
$ws.Range("F4").Value = "2021-10-05 13:38:29.769974"
$ws.Range("F5").Value = "2021-10-05 13:38:29.769976"
$ws.Range("F6").Value = "2021-10-05 13:38:29.769979"
$ws.Range("F7").Value = "2021-10-05 13:38:29.769982"
$ws.Range("F8").Value = "2021-10-05 13:38:29.769984"
$ws.Range("F9").Value = "2021-10-05 13:38:29.769987"
$ws.Range("F10").Value = "2021-10-05 13:38:29.769990"
$ws.Range("F11").Value = "2021-10-05 13:38:29.769992"
$ws.Range("F12").Value = "2021-10-05 13:38:29.769995"
$ws.Range("F13").Value = "2021-10-05 13:38:29.769997"
$ws.Range("F14").Value = "2021-10-05 13:38:29.770000"
